$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 32-39 (South Africa / Cape Town flight data continuing into December 2021) ---
# Copy formatting from the last existing data row (31) down through the new rows first,
# so the new cells inherit the same styles (s=1 for dates, s=2 for counts, s=4 for percent).
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D39").PasteSpecial(-4122)

$newRows = @(
    @{Row=32; Date="2021-11-30"; B=86;  C=56},
    @{Row=33; Date="2021-12-01"; B=103; C=60},
    @{Row=34; Date="2021-12-02"; B=104; C=64},
    @{Row=35; Date="2021-12-03"; B=114; C=71},
    @{Row=36; Date="2021-12-04"; B=82;  C=52},
    @{Row=37; Date="2021-12-05"; B=99;  C=63},
    @{Row=38; Date="2021-12-06"; B=99;  C=65},
    @{Row=39; Date="2021-12-07"; B=90;  C=56}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Formula = "=C$row/B$row"
}

# Clear the explicit row height so these rows use the default (matches how the source rows settled).
$ws.Range("A32:A39").EntireRow.AutoFit()

# --- Hyperlink to the CDC testing order (row 38, column F) ---
$ws.Hyperlinks.Add($ws.Range("F38"), "https://www.cdc.gov/coronavirus/2019-ncov/travelers/testing-international-air-travelers.html", "", "", "https://www.cdc.gov/coronavirus/2019-ncov/travelers/testing-international-air-travelers.html")
$ws.Range("F38").Value = "All flights departing after 12:01 a.m. ET December 6 will abide by a new CDC testing order."

# --- Footnote about the new CDC testing order (row 28, column F) ---
$ws.Range("F28").Value = "travel bans announced on November 26 bar entry into the US of noncitizens coming from eight countries in southern Africa. They are Botswana, Eswatini, Lesotho, Malawi, Mozambique, Namibia, South Africa and Zimbabwe. "
